$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Rotate the leaked OpenAI API key shown in the "OPENAI_API_KEY=" table
#    row, and collapse the now-superfluous empty paragraph that trailed it.
# ---------------------------------------------------------------------------
$oldKey = "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"
$newKey = "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"

$d.Content.Find.Execute($oldKey, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newKey, 2) | Out-Null

# Locate the table/cell that now holds the rotated key so we can drop the
# trailing blank paragraph that used to sit underneath it.
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $tbl = $d.Tables.Item($i)
    if ($tbl.Range.Text -like "*$newKey*") {
        for ($r = 1; $r -le $tbl.Rows.Count; $r++) {
            $row = $tbl.Rows.Item($r)
            for ($c = 1; $c -le $row.Cells.Count; $c++) {
                $cell = $row.Cells.Item($c)
                if ($cell.Range.Text -like "*$newKey*" -and $cell.Range.Paragraphs.Count -gt 1) {
                    $cell.Range.Paragraphs.Item(2).Range.Delete()
                }
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Tidy up the "3. Backend (Flask + OpenAI API)" label: it was split across
#    three runs around a spell-check proofErr wrapper on "OpenAI". Re-assert
#    the text in place so Word collapses it back into a single clean run.
# ---------------------------------------------------------------------------
$backendLabel = "3. Backend (Flask + OpenAI API)"
$d.Content.Find.Execute($backendLabel, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $backendLabel, 2) | Out-Null
